$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.977.78'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.326.10'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '187.99'
$ws.Range("E5").Value = '  +3.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '554.82'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.320.03'
$ws.Range("E8").Value = '  +2.24%  '
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.90'
$ws.Range("E12").Value = '  -2.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000264'
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.854.69'
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '578.87'
$ws.Range("E16").Value = '  -7.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.024.87'
$ws.Range("E17").Value = '  +1.10%  '
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.323.62'
$ws.Range("E19").Value = '  +1.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.77'
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.91'
$ws.Range("E21").Value = '  -3.52%  '
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.03'
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("E24").Value = '  +2.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.67'
$ws.Range("E25").Value = '  -6.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.96'
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("E27").Value = '  +1.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.36'
$ws.Range("E28").Value = '  -1.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '30.56'
$ws.Range("E29").Value = '  +1.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.41'
$ws.Range("E30").Value = '  -2.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.65'
$ws.Range("E31").Value = '  +6.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '574.67'
$ws.Range("E32").Value = '  +5.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.73'
$ws.Range("E33").Value = '  -4.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.85'
$ws.Range("E34").Value = '  -0.99%  '
$ws.Range("E35").Value = '  -1.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.718.72'
$ws.Range("E36").Value = '  +2.91%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.27'
$ws.Range("E38").Value = '  -2.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.00'
$ws.Range("E39").Value = '  +6.02%  '
$ws.Range("E40").Value = '  -2.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0688'
$ws.Range("E41").Value = '  -3.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.62'
$ws.Range("E42").Value = '  -3.28%  '
$ws.Range("E43").Value = '  -6.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.38'
$ws.Range("E44").Value = '  +3.13%  '
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0410'
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.128'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("B49").Value = 'CoreDAO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.90'
$ws.Range("E49").Value = '  -13.09%  '
$ws.Range("E50").Value = '  -3.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '126.75'
$ws.Range("E51").Value = '  +5.98%  '
